$d = $word.ActiveDocument

# 1. Remove the "lastRenderedPageBreak" rendering artifact on the paragraph that
#    starts with "La tabla "Pedidos"..." by touching its run with a no-op
#    find/replace (forces the run to be re-emitted without the stale marker).
$found = $d.Content.Find.Execute("La tabla", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "La tabla", 2)

# 2. Remove the second "ESTADO_PEDIDO((FK(id_estado)), (FK(id_pedido)))" paragraph
#    (near the end of the document) together with one of the two blank
#    paragraphs that precede the "EMPLEADO_DEPARTAMENTO" entry just above it.
$d.Paragraphs(84).Range.Delete()
$d.Paragraphs(81).Range.Delete()

# 3. Remove the first "ESTADO_PEDIDO((FK(id_ep)), (FK(id_estado)), fecha)" paragraph
#    and the empty "_GoBack" bookmark paragraph that immediately follows it.
$d.Paragraphs(32).Range.Delete()
$d.Paragraphs(31).Range.Delete()
